$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2023-12-01 Friday" "2023-12-02 Saturday"

Replace-Text "56×56=" "83×88="
Replace-Text "60×11=" "85×93="
Replace-Text "71×31=" "74×88="
Replace-Text "22×75=" "54×18="
Replace-Text "43×18=" "51×85="

Replace-Text "64×54=" "23×32="
Replace-Text "78×73=" "91×25="
Replace-Text "34×26=" "26×28="
# "26×22=" -> "13×21=" must run before "93×65=" -> "26×22=" below,
# otherwise the later replacement's output would collide with this source text.
Replace-Text "26×22=" "13×21="
Replace-Text "33×29=" "79×14="

Replace-Text "26×31=" "50×85="
Replace-Text "85×65=" "55×99="
Replace-Text "40×23=" "20×88="
Replace-Text "88×24=" "25×60="
Replace-Text "14×66=" "34×55="

Replace-Text "92×67=" "94×33="
Replace-Text "15×96=" "48×60="
Replace-Text "99×19=" "15×74="
Replace-Text "61×15=" "80×89="
Replace-Text "21×66=" "20×16="

Replace-Text "79×36=" "73×41="
Replace-Text "92×50=" "46×82="
Replace-Text "72×26=" "93×64="
Replace-Text "93×65=" "26×22="
Replace-Text "64×57=" "29×95="
